# Atualizado por script em 11-11-2023 20:45
#
# - Rows 3 and 4 (Petrovac-Mornar Bar / Decic-Jedinstvo) had their
#   home/away/odds/url data (columns F:V) swapped - index (A) and
#   kickoff date (E) stay put.
# - Rows 58 and 59 (Petrovac-Decic / Mornar Bar-Mladost DG) got the same
#   F:V swap treatment.
# - Four newly played fixtures were appended as rows 76-79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($sheet, $row1, $row2, $colStart, $colEnd) {
    $buffer = @{}
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $buffer[$c] = $sheet.Cells.Item($row1, $c).Value2
    }
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $sheet.Cells.Item($row1, $c).Value = $sheet.Cells.Item($row2, $c).Value2
    }
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $sheet.Cells.Item($row2, $c).Value = $buffer[$c]
    }
}

# --- Swap the two pairs of mixed-up fixtures (columns F..V = 6..22) ---
Swap-RowData $ws 3 4 6 22
Swap-RowData $ws 58 59 6 22

function Set-FixtureRow {
    param($sheet, $r, $idx, $e, $f, $g, $h, $i, $j, $k, $l, $m, $n, $o, $p, $q, $rr, $s, $t, $u, $v)

    $sheet.Cells.Item($r, 1).Value = $idx
    $sheet.Cells.Item(75, 1).Copy()
    $sheet.Cells.Item($r, 1).PasteSpecial(-4122)

    $sheet.Cells.Item($r, 2).Value = "montenegro"
    $sheet.Cells.Item($r, 3).Value = "prva-crnogorska-liga"
    $sheet.Cells.Item($r, 4).Value = "2023-2024"

    $sheet.Cells.Item($r, 5).Value = $e
    $sheet.Cells.Item(75, 5).Copy()
    $sheet.Cells.Item($r, 5).PasteSpecial(-4122)

    $sheet.Cells.Item($r, 6).Value = $f
    $sheet.Cells.Item($r, 7).Value = $g
    $sheet.Cells.Item($r, 8).Value = $h
    $sheet.Cells.Item($r, 9).Value = $i
    $sheet.Cells.Item($r, 10).Value = $j
    $sheet.Cells.Item($r, 11).Value = $k
    $sheet.Cells.Item($r, 12).Value = $l
    $sheet.Cells.Item($r, 13).Value = $m
    $sheet.Cells.Item($r, 14).Value = $n
    $sheet.Cells.Item($r, 15).Value = $o
    $sheet.Cells.Item($r, 16).Value = $p
    $sheet.Cells.Item($r, 17).Value = $q
    $sheet.Cells.Item($r, 18).Value = $rr
    $sheet.Cells.Item($r, 19).Value = $s
    $sheet.Cells.Item($r, 20).Value = $t
    $sheet.Cells.Item($r, 21).Value = $u
    $sheet.Cells.Item($r, 22).Value = $v
}

# --- Append the four new fixtures that were scraped on 11/11/2023 ---

Set-FixtureRow $ws 76 75 45241.66666666666 `
    "Mladost DG" 1 "Sutjeska" 1 `
    4.61 "10/11/2023 04:13" 6.15 "11/11/2023 15:43" `
    3.5 "10/11/2023 04:13" 3.73 "11/11/2023 15:43" `
    1.63 "10/11/2023 04:13" 1.56 "11/11/2023 15:43" `
    "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/mladost-dg-sutjeska/G6X7cGLN/"

Set-FixtureRow $ws 77 76 45241.66666666666 `
    "Buducnost" 3 "Jedinstvo" 2 `
    1.29 "10/11/2023 04:13" 1.29 "11/11/2023 15:42" `
    4.83 "10/11/2023 04:13" 5.32 "11/11/2023 15:45" `
    7.5 "10/11/2023 04:13" 8.96 "11/11/2023 15:42" `
    "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/buducnost-jedinstvo/0v9vhfqo/"

Set-FixtureRow $ws 78 77 45241.66666666666 `
    "Petrovac" 1 "Arsenal Tivat" 1 `
    2.05 "10/11/2023 04:13" 1.97 "11/11/2023 15:43" `
    2.98 "10/11/2023 04:13" 3.06 "11/11/2023 15:43" `
    3.43 "10/11/2023 04:13" 4.21 "11/11/2023 15:43" `
    "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/petrovac-arsenal-tivat/bVMCdzyU/"

Set-FixtureRow $ws 79 78 45241.77083333334 `
    "Decic" 0 "Jezero" 0 `
    1.58 "10/11/2023 06:42" 1.5 "11/11/2023 18:28" `
    3.44 "10/11/2023 06:42" 3.63 "11/11/2023 18:28" `
    5.22 "10/11/2023 06:42" 7.92 "11/11/2023 18:28" `
    "https://www.betexplorer.com/football/montenegro/prva-crnogorska-liga/decic-jezero/23T3bd6H/"
